$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update column F (dSF) values to match repulled data
$ws.Range("F2").Value = 12
$ws.Range("F3").Value = -2
$ws.Range("F6").Value = -3
$ws.Range("F9").Value = -2
$ws.Range("F14").Value = -8
